$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header text (MODEL_CONDITION -> MODELCONDITION) before
# the column shift, while it still lives in E1.
$ws.Range("E1").Value = "MODELCONDITION"

# The old column A (style-only "11"/"17" helper column) is removed; every
# other column shifts one place to the left (B->A, C->B, D->C, E->D, F->E).
$ws.Columns.Item(1).Delete()
